$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Controllers that handle dips in pH run longer: bump the "Step" duration
# for each configured controller (rows 2-6, column D) from 60 to 120.
$ws.Range("D2:D6").Value = 120

# Leave the selection on the range that was just edited.
$ws.Range("D2:D6").Select() | Out-Null
